# The deck's "datetimeFigureOut" date placeholders (on the slide master and
# every slide layout) were refreshed by PowerPoint from 9/27/2022 to
# 6/12/2023 — the ordinary side effect of PowerPoint recalculating the
# auto-updating date field the next time the file was opened/saved, rather
# than a deliberate content edit on either slide.
#
# ppPlaceholderDate == 16 identifies the date placeholder shape regardless
# of its index/name on a given master or layout.

$p = $ppt.ActivePresentation
$newDate = "6/12/2023"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own date placeholder.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li).Shapes
}
